# Update odds/score-distribution figures for the 2025-02-05 FlashScore
# "Jogos da Semana" sheet. Only numeric odds/count cells change; the
# match metadata (Id/Date/Time/League/Home/Away) in columns A-F is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Zalaegerszegi vs Ferencvaros
$ws.Range("H3").Value  = 3.75
$ws.Range("J3").Value  = 5.5
$ws.Range("K3").Value  = 2.25
$ws.Range("L3").Value  = 2.2
$ws.Range("Q3").Value  = 1.9
$ws.Range("R3").Value  = 1.95
$ws.Range("S3").Value  = 3.25
$ws.Range("T3").Value  = 1.33
$ws.Range("AB3").Value = 51
$ws.Range("AE3").Value = 10
$ws.Range("AL3").Value = 12
$ws.Range("AO3").Value = 301

# Row 4 - Mohun Bagan vs Punjab
$ws.Range("G4").Value  = 1.48
$ws.Range("H4").Value  = 4.1
$ws.Range("I4").Value  = 5.75
$ws.Range("K4").Value  = 2.4
$ws.Range("M4").Value  = 1.04
$ws.Range("N4").Value  = 13
$ws.Range("O4").Value  = 1.2
$ws.Range("P4").Value  = 4.33
$ws.Range("Q4").Value  = 1.65
$ws.Range("R4").Value  = 2.2
$ws.Range("S4").Value  = 2.63
$ws.Range("T4").Value  = 1.44
$ws.Range("W4").Value  = 1.8
$ws.Range("X4").Value  = 1.91
$ws.Range("Y4").Value  = 8
$ws.Range("Z4").Value  = 7.5
$ws.Range("AE4").Value = 13
$ws.Range("AP4").Value = 2.03
$ws.Range("AQ4").Value = 1.78

# Row 5 - Johor DT vs PDRM FC
$ws.Range("G5").Value  = 1.03
$ws.Range("H5").Value  = 8.5
$ws.Range("I5").Value  = 40
$ws.Range("J5").Value  = 1.23
$ws.Range("K5").Value  = 3.6
$ws.Range("L5").Value  = 28
$ws.Range("Q5").Value  = 1.23
$ws.Range("R5").Value  = 3.9
$ws.Range("S5").Value  = 1.55
$ws.Range("T5").Value  = 2.15
$ws.Range("Y5").Value  = 10
$ws.Range("Z5").Value  = 5.9
$ws.Range("AA5").Value = 15
$ws.Range("AB5").Value = 4.9
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 55
$ws.Range("AE5").Value = 20
$ws.Range("AF5").Value = 27
$ws.Range("AG5").Value = 70
$ws.Range("AH5").Value = 400
$ws.Range("AI5").Value = 175
$ws.Range("AJ5").Value = 500
$ws.Range("AK5").Value = 250

# Row 6 - Pachuca vs Club Leon
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2

# Row 8 - St. Gallen vs Lugano
$ws.Range("G8").Value  = 2.7
$ws.Range("I8").Value  = 2.5
$ws.Range("J8").Value  = 3.2
$ws.Range("L8").Value  = 3.1
$ws.Range("Q8").Value  = 1.7
$ws.Range("R8").Value  = 2.1
$ws.Range("Z8").Value  = 15
$ws.Range("AC8").Value = 21
$ws.Range("AD8").Value = 26
$ws.Range("AJ8").Value = 13
$ws.Range("AM8").Value = 19
$ws.Range("AN8").Value = 23
